$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header: Q8 in J1 (copy formatting from the preceding header cell
# so it matches the bold/bordered/centered header style).
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null
$ws.Range("J1").Value = "Q8"

# Full refresh of the data block (values were regenerated upstream, and one
# extra column of simulated quarters was added for several rows).
$values = @{
    2  = @(-0.740857461610962, 0.2483496536967165, 0.8032852183307098, 1.015296315185831, 0.4272448182495295, 0.01108471313272752, 1.401227768176947)
    3  = @(0.2917404373296907, 0.8466760019636841, 1.058687098818806, 0.4706356018825037, 0.05447549676570174, 1.444618551809921)
    4  = @(1.247683574918372, 1.459694671773494, 0.8716431748371921, 0.4554830697203902, 1.84562612476461, 1.076491995083831, -0.7442442573846902, 1.069262038377534, 0.4350014876132097)
    5  = @(3.694617372890321, 3.106565875954019, 2.690405770837217, 4.080548825881436, 3.311414696200658, 1.490678443732137, 3.304184739494361, 2.669924188730037)
    6  = @(1.216530487278416, 0.8003703821616144, 2.190513437205834, 1.421379307525055, -0.399356944943466, 1.414149350818758, 0.779888800054434)
    7  = @(0.4476790584865185, 1.837822113530738, 1.068687983849959, -0.7520482686185619, 1.061458027143662, 0.427197476379338)
    8  = @(1.594404170131267, 0.8252700404504878, -0.9954662120180333, 0.8180400837441908, 0.1837795329798666, 1.265495818666463, 0.8840541853673727, 0.5915919440004813)
    9  = @(0.4048306212132332, -1.415905631255288, 0.3976006645069362, -0.236659886257388, 0.8450563994292083, 0.4636147661301181, 0.1711525247632267)
    10 = @(-1.433992460878194, 0.3795138348840296, -0.2547467158802946, 0.8269695698063018, 0.4455279365072115, 0.1530656951403201)
    11 = @(0.4299722955860048, -0.2042882551783194, 0.8774280305082769, 0.4959863972091867, 0.2035241558422953)
    12 = @(-0.4160968922281114, 0.6656193934584849, 0.2841777601593947, -0.008284481207496679)
    13 = @(0.5354267536149976, 0.1539851203159074, -0.1384771210509839)
    14 = @(-0.322788625881465, -0.6152508672483563)
    15 = @(-0.2720993704486361)
    16 = @()
}

# Clear any stale values beyond the new row extents (old sheet had up to
# column I with different per-row lengths); overwrite B:J for rows 2-16.
foreach ($r in 2..16) {
    foreach ($c in 2..10) {
        $ws.Cells.Item($r, $c).Value = $null
    }
}

foreach ($r in $values.Keys) {
    $row = $values[$r]
    for ($i = 0; $i -lt $row.Count; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($r, $col).Value = $row[$i]
    }
}
